$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation is inserted as row 9; every existing
# record from the old row 9 down to the old last row (43) shifts down by
# one row (to 10..44).
$ws.Rows("9:9").Insert()

$ws.Cells.Item(9, 1).Value = 5
$ws.Cells.Item(9, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(9, 3).Value = "Maule"
$ws.Cells.Item(9, 4).Value = Get-Date -Year 2022 -Month 5 -Day 9 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(9, 5).Value = 7
$ws.Cells.Item(9, 6).Value = "Fruta"
$ws.Cells.Item(9, 7).Value = 100104
$ws.Cells.Item(9, 8).Value = "Frutos de pepita"
$ws.Cells.Item(9, 9).Value = 100104003
$ws.Cells.Item(9, 10).Value = "Membrillo"
$ws.Cells.Item(9, 11).Value = "Champion"
$ws.Cells.Item(9, 12).Value = "Primera"
$ws.Cells.Item(9, 13).Value = 230
$ws.Cells.Item(9, 14).Value = 10000
$ws.Cells.Item(9, 15).Value = 10000
$ws.Cells.Item(9, 16).Value = 10000
$ws.Cells.Item(9, 17).Value = '$/caja 18 kilos granel'
$ws.Cells.Item(9, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(9, 19).Value = 556
$ws.Cells.Item(9, 20).Value = 18
